$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.005243333333333
$ws.Range("H2").Value = 3.01573
$ws.Range("I2").Value = 0.07224874268505826
$ws.Range("J2").Value = 0.07224874268505825
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.616180666666666
$ws.Range("N2").Value = 22.848542
$ws.Range("O2").Value = 0.06266940401417194
$ws.Range("P2").Value = 0.06266940401417194
$ws.Range("Q2").Value = 7.656114840628888
$ws.Range("R2").Value = 68.90503356565999
$ws.Range("S2").Value = 0.004527785644845866
$ws.Range("T2").Value = 0.004527785644845865

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.005243333333333
$ws.Range("H3").Value = 3.01573
$ws.Range("I3").Value = 0.07224874268505826
$ws.Range("J3").Value = 0.07224874268505825
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 86.98680366666667
$ws.Range("N3").Value = 260.960411
$ws.Range("O3").Value = 0.7157670466966058
$ws.Range("P3").Value = 0.7157670466966058
$ws.Range("Q3").Value = 87.44290447389223
$ws.Range("R3").Value = 786.9861402650299
$ws.Range("S3").Value = 0.05171326917922716
$ws.Range("T3").Value = 0.05171326917922715

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.005243333333333
$ws.Range("H4").Value = 3.01573
$ws.Range("I4").Value = 0.07224874268505826
$ws.Range("J4").Value = 0.07224874268505825
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.9360126666666666
$ws.Range("N4").Value = 2.808038
$ws.Range("O4").Value = 0.007701938614251506
$ws.Range("P4").Value = 0.007701938614251506
$ws.Range("Q4").Value = 0.9409204930822221
$ws.Range("R4").Value = 8.468284437739998
$ws.Range("S4").Value = 0.0005564553811171712
$ws.Range("T4").Value = 0.0005564553811171711

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.005243333333333
$ws.Range("H5").Value = 3.01573
$ws.Range("I5").Value = 0.07224874268505826
$ws.Range("J5").Value = 0.07224874268505825
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 25.99049233333334
$ws.Range("N5").Value = 77.97147700000001
$ws.Range("O5").Value = 0.2138616106749707
$ws.Range("P5").Value = 0.2138616106749707
$ws.Range("Q5").Value = 26.12676914813445
$ws.Range("R5").Value = 235.14092233321
$ws.Range("S5").Value = 0.01545123247986807
$ws.Range("T5").Value = 0.01545123247986806

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.25983933333333
$ws.Range("H6").Value = 30.779518
$ws.Range("I6").Value = 0.7373940889775011
$ws.Range("J6").Value = 0.737394088977501
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.616180666666666
$ws.Range("N6").Value = 22.848542
$ws.Range("O6").Value = 0.06266940401417194
$ws.Range("P6").Value = 0.06266940401417194
$ws.Range("Q6").Value = 78.14078997363956
$ws.Range("R6").Value = 703.2671097627559
$ws.Range("S6").Value = 0.04621204807979327
$ws.Range("T6").Value = 0.04621204807979326

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.25983933333333
$ws.Range("H7").Value = 30.779518
$ws.Range("I7").Value = 0.7373940889775011
$ws.Range("J7").Value = 0.737394088977501
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 86.98680366666667
$ws.Range("N7").Value = 260.960411
$ws.Range("O7").Value = 0.7157670466966058
$ws.Range("P7").Value = 0.7157670466966058
$ws.Range("Q7").Value = 892.470629740211
$ws.Range("R7").Value = 8032.235667661898
$ws.Range("S7").Value = 0.5278023893189602
$ws.Range("T7").Value = 0.5278023893189601

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.25983933333333
$ws.Range("H8").Value = 30.779518
$ws.Range("I8").Value = 0.7373940889775011
$ws.Range("J8").Value = 0.737394088977501
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.9360126666666666
$ws.Range("N8").Value = 2.808038
$ws.Range("O8").Value = 0.007701938614251506
$ws.Range("P8").Value = 0.007701938614251506
$ws.Range("Q8").Value = 9.603339573964888
$ws.Range("R8").Value = 86.430056165684
$ws.Range("S8").Value = 0.005679364007816627
$ws.Range("T8").Value = 0.005679364007816626

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.25983933333333
$ws.Range("H9").Value = 30.779518
$ws.Range("I9").Value = 0.7373940889775011
$ws.Range("J9").Value = 0.737394088977501
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 25.99049233333334
$ws.Range("N9").Value = 77.97147700000001
$ws.Range("O9").Value = 0.2138616106749707
$ws.Range("P9").Value = 0.2138616106749707
$ws.Range("Q9").Value = 266.6582755342318
$ws.Range("R9").Value = 2399.924479808086
$ws.Range("S9").Value = 0.157700287570931
$ws.Range("T9").Value = 0.157700287570931

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.383875
$ws.Range("H10").Value = 1.151625
$ws.Range("I10").Value = 0.02758982345723265
$ws.Range("J10").Value = 0.02758982345723265
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.616180666666666
$ws.Range("N10").Value = 22.848542
$ws.Range("O10").Value = 0.06266940401417194
$ws.Range("P10").Value = 0.06266940401417194
$ws.Range("Q10").Value = 2.923661353416667
$ws.Range("R10").Value = 26.31295218075
$ws.Range("S10").Value = 0.001729037792920991
$ws.Range("T10").Value = 0.001729037792920991

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.383875
$ws.Range("H11").Value = 1.151625
$ws.Range("I11").Value = 0.02758982345723265
$ws.Range("J11").Value = 0.02758982345723265
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 86.98680366666667
$ws.Range("N11").Value = 260.960411
$ws.Range("O11").Value = 0.7157670466966058
$ws.Range("P11").Value = 0.7157670466966058
$ws.Range("Q11").Value = 33.39205925754167
$ws.Range("R11").Value = 300.528533317875
$ws.Range("S11").Value = 0.01974788645486416
$ws.Range("T11").Value = 0.01974788645486416

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.383875
$ws.Range("H12").Value = 1.151625
$ws.Range("I12").Value = 0.02758982345723265
$ws.Range("J12").Value = 0.02758982345723265
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.9360126666666666
$ws.Range("N12").Value = 2.808038
$ws.Range("O12").Value = 0.007701938614251506
$ws.Range("P12").Value = 0.007701938614251506
$ws.Range("Q12").Value = 0.3593118624166667
$ws.Range("R12").Value = 3.23380676175
$ws.Range("S12").Value = 0.0002124951266456422
$ws.Range("T12").Value = 0.0002124951266456422

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.383875
$ws.Range("H13").Value = 1.151625
$ws.Range("I13").Value = 0.02758982345723265
$ws.Range("J13").Value = 0.02758982345723265
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 25.99049233333334
$ws.Range("N13").Value = 77.97147700000001
$ws.Range("O13").Value = 0.2138616106749707
$ws.Range("P13").Value = 0.2138616106749707
$ws.Range("Q13").Value = 9.977100244458335
$ws.Range("R13").Value = 89.79390220012502
$ws.Range("S13").Value = 0.005900404082801863
$ws.Range("T13").Value = 0.005900404082801863

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.264687
$ws.Range("H14").Value = 6.794061
$ws.Range("I14").Value = 0.162767344880208
$ws.Range("J14").Value = 0.162767344880208
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 7.616180666666666
$ws.Range("N14").Value = 22.848542
$ws.Range("O14").Value = 0.06266940401417194
$ws.Range("P14").Value = 0.06266940401417194
$ws.Range("Q14").Value = 17.24826534545133
$ws.Range("R14").Value = 155.234388109062
$ws.Range("S14").Value = 0.01020053249661182
$ws.Range("T14").Value = 0.01020053249661182

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.264687
$ws.Range("H15").Value = 6.794061
$ws.Range("I15").Value = 0.162767344880208
$ws.Range("J15").Value = 0.162767344880208
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 86.98680366666667
$ws.Range("N15").Value = 260.960411
$ws.Range("O15").Value = 0.7157670466966058
$ws.Range("P15").Value = 0.7157670466966058
$ws.Range("Q15").Value = 196.9978834354523
$ws.Range("R15").Value = 1772.980950919071
$ws.Range("S15").Value = 0.1165035017435544
$ws.Range("T15").Value = 0.1165035017435544

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2.264687
$ws.Range("H16").Value = 6.794061
$ws.Range("I16").Value = 0.162767344880208
$ws.Range("J16").Value = 0.162767344880208
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.9360126666666666
$ws.Range("N16").Value = 2.808038
$ws.Range("O16").Value = 0.007701938614251506
$ws.Range("P16").Value = 0.007701938614251506
$ws.Range("Q16").Value = 2.119775718035333
$ws.Range("R16").Value = 19.077981462318
$ws.Range("S16").Value = 0.001253624098672066
$ws.Range("T16").Value = 0.001253624098672066

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 2.264687
$ws.Range("H17").Value = 6.794061
$ws.Range("I17").Value = 0.162767344880208
$ws.Range("J17").Value = 0.162767344880208
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 25.99049233333334
$ws.Range("N17").Value = 77.97147700000001
$ws.Range("O17").Value = 0.2138616106749707
$ws.Range("P17").Value = 0.2138616106749707
$ws.Range("Q17").Value = 58.86033011089967
$ws.Range("R17").Value = 529.742970998097
$ws.Range("S17").Value = 0.03480968654136973
$ws.Range("T17").Value = 0.03480968654136973
